$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (line7 replaces what used to be extr1 data in this row)
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9 (line8 replaces what used to be extr2 data in this row)
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# Row 10 (extr1, shifted down from previous row 8 data position)
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (extr2)
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12 (extr3)
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# Row 13 (extr4)
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14 (extr5)
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# Row 15 (extr6)
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# Row 16 (extr7) - new row
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 1).Font.Bold = $true
$ws.Cells.Item(16, 1).HorizontalAlignment = -4108
$ws.Cells.Item(16, 1).VerticalAlignment = -4160
$ws.Cells.Item(16, 1).Borders.LineStyle = 1
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# Row 17 (extr8) - new row
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 1).Font.Bold = $true
$ws.Cells.Item(17, 1).HorizontalAlignment = -4108
$ws.Cells.Item(17, 1).VerticalAlignment = -4160
$ws.Cells.Item(17, 1).Borders.LineStyle = 1
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
